# The document had the tag content "<id>p017r_1</id>" split across three
# separate runs (one for "<id>", one for "p017r_1", one for "</id>").
# The edit merges them into a single run containing the full text
# "<id>p017r_1</id>".
#
# Doing a Find/Replace over the exact phrase spanning all three runs
# collapses the matched range into one run, taking on the character
# formatting of the first run in the match (Courier New / color 7f6000 /
# sz 18), which is exactly the formatting the merged run keeps in the
# target diff.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p017r_1</id>",   # FindText
    $true,                # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "<id>p017r_1</id>",   # ReplaceWith
    2                     # Replace (wdReplaceAll)
)
